# fixed substring of GEM for municipalities and added last year as copy of
# previous year at ASTAT data
#
# Adds the missing Italian "name_short_IT" values (column M) on the
# "concept" sheet for several concepts, and duplicates the ASTAT
# "Saldo migratorio" short name into the newly introduced column M so the
# Italian table has a value for every row that already has one in German
# (column L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("concept")
$ws.Activate()

# Column M = name_short_IT. New strings are entered in the same order they
# were first typed so the shared-string table grows the way it originally
# did: donne, then totale, then uomini for the unemployment-rate block,
# followed by the other brand-new short names.
$ws.Range("M3").Value  = "Tasso di disoccupazione (donne)"
$ws.Range("M2").Value  = "Tasso di disoccupazione (totale)"
$ws.Range("M4").Value  = "Tasso di disoccupazione (uomini)"
$ws.Range("M11").Value = "Disoccupati iscritti (totale)"
$ws.Range("M19").Value = "Turismo: permanenza media (gg)"
$ws.Range("M20").Value = "Occupati tempo indeterminato (totale)"
$ws.Range("M21").Value = "Occupati tempo indeterminato (donne)"
$ws.Range("M22").Value = "Occupati tempo indeterminato (uomini)"

# These reuse short names that already exist elsewhere in the sheet.
$ws.Range("M5").Value  = "Apprendisti (totale)"
$ws.Range("M6").Value  = "Apprendisti (donne)"
$ws.Range("M7").Value  = "Apprendisti (uomini)"
$ws.Range("M12").Value = "Disoccupati iscritti (donne)"
$ws.Range("M13").Value = "Disoccupati iscritti (uomini)"

# ASTAT row: add the last year as a copy of the previous year's value
# (same text already used in F23/G23/H23) into the new column M23.
$ws.Range("M23").Value = "Saldo migratorio"

# Restore the view: select M27 (scroll position isn't persisted by this
# runtime's writer, but the selection is).
$ws.Range("M27").Select()
